$wb = $excel.ActiveWorkbook

# --- "Volume Ranges" sheet: add two derived columns -----------------------
# D = Voltage Amplitude (Speaker Volume = 42)  -> half of column B
# E = Voltage amplitdue (Speaker Volume = 100) -> half of column C
$wsVol = $wb.Worksheets.Item("Volume Ranges")

$wsVol.Range("D1").Value = "Voltage Amplitude (Speaker Volume = 42)"
$wsVol.Range("E1").Value = "Voltage amplitdue (Speaker Volume = 100)"

# Row 2 gets its own (non-shared) formula, rows 3-12 fill down as a block
# so they share one formula group - matches how Excel would record a
# single cell entered first, then the rest filled down together.
$wsVol.Range("D2").Formula = "=B2/2"
$wsVol.Range("E2").Formula = "=C2/2"
$wsVol.Range("D3:D12").Formula = "=B3/2"
$wsVol.Range("E3:E12").Formula = "=C3/2"

# Size the two new columns to fit their (longer) header text, same as the
# existing B:C columns.
$wsVol.Columns.Item(4).AutoFit() | Out-Null
$wsVol.Columns.Item(5).AutoFit() | Out-Null

# --- "Percentage Division" sheet: drop the stray italic-ish formatting ----
# that had been applied (then undone) to rows 4, 8 and 14 of A:C.
$wsPct = $wb.Worksheets.Item("Percentage Division")
$wsPct.Range("A4:C4").Style = "Normal"
$wsPct.Range("A8:C8").Style = "Normal"
$wsPct.Range("A14:C14").Style = "Normal"

# --- Active sheet / selection bookkeeping ----------------------------------
# "Volume Ranges" becomes the active tab, with D2:D12 having just been
# filled, the new E2:E12 column is left selected.
$wsVol.Activate() | Out-Null
$wsVol.Range("E2:E12").Select() | Out-Null

Write-Output "edit applied"
